{"js": "const replacements = [\n  [\"36\u00f74=9, 0\", \"37\u00f73=12, 1\"],\n  [\"15\u00f77=2, 1\", \"30\u00f79=3, 3\"],\n  [\"86\u00f77=12, 2\", \"92\u00f77=13, 1\"],\n  [\"81\u00f75=16, 1\", \"82\u00f77=11, 5\"],\n  [\"40\u00f79=4, 4\", \"81\u00f76=13, 3\"],\n  [\"81\u00f78=10, 1\", \"99\u00f78=12, 3\"],\n  [\"11\u00f72=5, 1\", \"51\u00f76=8, 3\"],\n  [\"32\u00f74=8, 0\", \"34\u00f77=4, 6\"],\n  [\"12\u00f79=1, 3\", \"66\u00f77=9, 3\"],\n  [\"97\u00f78=12, 1\", \"65\u00f78=8, 1\"],\n  [\"14\u00f74=3, 2\", \"34\u00f76=5, 4\"],\n  [\"66\u00f72=33, 0\", \"45\u00f78=5, 5\"],\n  [\"15\u00f75=3, 0\", \"85\u00f74=21, 1\"],\n  [\"64\u00f76=10, 4\", \"33\u00f78=4, 1\"],\n  [\"14\u00f79=1, 5\", \"29\u00f79=3, 2\"],\n  [\"97\u00f75=19, 2\", \"70\u00f75=14, 0\"],\n  [\"39\u00f79=4, 3\", \"35\u00f76=5, 5\"],\n  [\"42\u00f76=7, 0\", \"31\u00f75=6, 1\"],\n  [\"22\u00f73=7, 1\", \"63\u00f74=15, 3\"],\n  [\"50\u00f78=6, 2\", \"13\u00f76=2, 1\"],\n  [\"88\u00f74=22, 0\", \"95\u00f78=11, 7\"],\n  [\"59\u00f76=9, 5\", \"83\u00f72=41, 1\"],\n  [\"51\u00f75=10, 1\", \"85\u00f76=14, 1\"],\n  [\"88\u00f79=9, 7\", \"44\u00f79=4, 8\"],\n  [\"65\u00f72=32, 1\", \"96\u00f76=16, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"36\u00f74=9, 0\", \"37\u00f73=12, 1\"),\n    @(\"15\u00f77=2, 1\", \"30\u00f79=3, 3\"),\n    @(\"86\u00f77=12, 2\", \"92\u00f77=13, 1\"),\n    @(\"81\u00f75=16, 1\", \"82\u00f77=11, 5\"),\n    @(\"40\u00f79=4, 4\", \"81\u00f76=13, 3\"),\n    @(\"81\u00f78=10, 1\", \"99\u00f78=12, 3\"),\n    @(\"11\u00f72=5, 1\", \"51\u00f76=8, 3\"),\n    @(\"32\u00f74=8, 0\", \"34\u00f77=4, 6\"),\n    @(\"12\u00f79=1, 3\", \"66\u00f77=9, 3\"),\n    @(\"97\u00f78=12, 1\", \"65\u00f78=8, 1\"),\n    @(\"14\u00f74=3, 2\", \"34\u00f76=5, 4\"),\n    @(\"66\u00f72=33, 0\", \"45\u00f78=5, 5\"),\n    @(\"15\u00f75=3, 0\", \"85\u00f74=21, 1\"),\n    @(\"64\u00f76=10, 4\", \"33\u00f78=4, 1\"),\n    @(\"14\u00f79=1, 5\", \"29\u00f79=3, 2\"),\n    @(\"97\u00f75=19, 2\", \"70\u00f75=14, 0\"),\n    @(\"39\u00f79=4, 3\", \"35\u00f76=5, 5\"),\n    @(\"42\u00f76=7, 0\", \"31\u00f75=6, 1\"),\n    @(\"22\u00f73=7, 1\", \"63\u00f74=15, 3\"),\n    @(\"50\u00f78=6, 2\", \"13\u00f76=2, 1\"),\n    @(\"88\u00f74=22, 0\", \"95\u00f78=11, 7\"),\n    @(\"59\u00f76=9, 5\", \"83\u00f72=41, 1\"),\n    @(\"51\u00f75=10, 1\", \"85\u00f76=14, 1\"),\n    @(\"88\u00f79=9, 7\", \"44\u00f79=4, 8\"),\n    @(\"65\u00f72=32, 1\", \"96\u00f76=16, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
